$wb = $excel.ActiveWorkbook

# Sheet: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 373
$ws1.Range("F6").Value = 768
$ws1.Range("F8").Value = 1082
$ws1.Range("F9").Value = 279
$ws1.Range("F12").Value = 629
$ws1.Range("F18").Value = 832
$ws1.Range("F31").Value = 240
$ws1.Range("F32").Value = 1033

# Sheet: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 1048
$ws2.Range("F5").Value = 1048
$ws2.Range("F22").Value = 35
$ws2.Range("F26").Value = 3746

# Sheet: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 1769
$ws3.Range("F10").Value = 332

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1769
$ws4.Range("F9").Value = 332
$ws4.Range("F12").Value = 373
$ws4.Range("F13").Value = 768
$ws4.Range("F16").Value = 1082
$ws4.Range("F17").Value = 279
$ws4.Range("F19").Value = 629
$ws4.Range("F20").Value = 1048
$ws4.Range("F24").Value = 832
$ws4.Range("F38").Value = 240
$ws4.Range("F42").Value = 35
$ws4.Range("F46").Value = 1033
